$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 89: Isaacson (2023) - insomnia - 3-13 weeks
$ws.Range("A89").Value = "Isaacson (2023)"
$ws.Range("B89").Value = "insomnia"
$ws.Range("C89").Value = "3-13 weeks"
$ws.Range("D89").Value = "Low"
$ws.Range("E89").Value = "Low"
$ws.Range("F89").Value = "Some concerns"
$ws.Range("G89").Value = "Low"
$ws.Range("H89").Value = "Low"
$ws.Range("J89").Value = "Low"

# Row 90: Tsukada (2023) - sedation - 1 day-2 weeks
$ws.Range("A90").Value = "Tsukada (2023)"
$ws.Range("B90").Value = "sedation"
$ws.Range("C90").Value = "1 day-2 weeks"
$ws.Range("D90").Value = "Some concerns"
$ws.Range("E90").Value = "Low"
$ws.Range("F90").Value = "Some concerns"
$ws.Range("G90").Value = "Low"
$ws.Range("H90").Value = "Low"
$ws.Range("I90").Value = "Some concerns"
$ws.Range("J90").Value = "Some concerns"

# Copy row styles from nearby existing rows (88 uses the same outcome/timepoint formatting pattern)
$ws.Range("A88:J88").Copy() | Out-Null
$ws.Range("A90:J90").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A26:J26").Copy() | Out-Null
$ws.Range("A89:J89").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A89").RowHeight = 34
$ws.Range("A90").RowHeight = 34

# Scroll the view down towards the new rows and select the last-edited cell,
# matching where the author ended up after adding the two rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C90").Select() | Out-Null
